$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'23.946.57"
$ws.Cells.Item(2, 5).Value = '  -1.90%  '
$ws.Cells.Item(3, 4).Value = "'1.623.96"
$ws.Cells.Item(3, 5).Value = '  -1.75%  '
$ws.Cells.Item(4, 4).Value = "'0.9972"
$ws.Cells.Item(4, 5).Value = '  -0.61%  '
$ws.Cells.Item(5, 4).Value = "'307.72"
$ws.Cells.Item(5, 5).Value = '  -1.33%  '
$ws.Cells.Item(6, 4).Value = "'0.9985"
$ws.Cells.Item(6, 5).Value = '  -0.55%  '
$ws.Cells.Item(7, 4).Value = "'0.3936"
$ws.Cells.Item(7, 5).Value = '  +0.59%  '
$ws.Cells.Item(8, 4).Value = "'0.3852"
$ws.Cells.Item(8, 5).Value = '  -1.38%  '
$ws.Cells.Item(9, 4).Value = "'0.9971"
$ws.Cells.Item(9, 5).Value = '  -0.56%  '
$ws.Cells.Item(10, 4).Value = "'49.53"
$ws.Cells.Item(10, 5).Value = '  -1.86%  '
$ws.Cells.Item(11, 4).Value = "'1.358"
$ws.Cells.Item(11, 5).Value = '  -1.34%  '
$ws.Cells.Item(12, 4).Value = "'0.08485"
$ws.Cells.Item(12, 5).Value = '  -0.91%  '
$ws.Cells.Item(13, 4).Value = "'23.84"
$ws.Cells.Item(13, 5).Value = '  -4.73%  '
$ws.Cells.Item(14, 4).Value = "'7.085"
$ws.Cells.Item(14, 5).Value = '  -1.61%  '
$ws.Cells.Item(15, 4).Value = "'7.626"
$ws.Cells.Item(15, 5).Value = '  +0.16%  '
$ws.Cells.Item(16, 4).Value = "'0.00001286"
$ws.Cells.Item(16, 5).Value = '  -1.42%  '
$ws.Cells.Item(17, 4).Value = "'1.622.03"
$ws.Cells.Item(17, 5).Value = '  -1.97%  '
$ws.Cells.Item(18, 4).Value = "'93.92"
$ws.Cells.Item(18, 5).Value = '  +0.82%  '
$ws.Cells.Item(19, 4).Value = "'0.06927"
$ws.Cells.Item(19, 5).Value = '  -0.35%  '
$ws.Cells.Item(20, 4).Value = "'20.04"
$ws.Cells.Item(20, 5).Value = '  -5.17%  '
$ws.Cells.Item(21, 4).Value = "'6.845"
$ws.Cells.Item(21, 5).Value = '  -2.23%  '
$ws.Cells.Item(22, 4).Value = "'0.9987"
$ws.Cells.Item(22, 5).Value = '  -0.61%  '
$ws.Cells.Item(23, 4).Value = "'13.44"
$ws.Cells.Item(23, 5).Value = '  -2.56%  '
$ws.Cells.Item(24, 4).Value = "'23.931.48"
$ws.Cells.Item(24, 5).Value = '  -1.96%  '
$ws.Cells.Item(25, 4).Value = "'2.477"
$ws.Cells.Item(25, 5).Value = '  +5.43%  '
$ws.Cells.Item(26, 4).Value = "'2.840"
$ws.Cells.Item(26, 5).Value = '  +2.25%  '
$ws.Cells.Item(27, 4).Value = "'22.23"
$ws.Cells.Item(27, 5).Value = '  -1.86%  '
$ws.Cells.Item(28, 4).Value = "'156.85"
$ws.Cells.Item(28, 5).Value = '  -1.11%  '
$ws.Cells.Item(29, 4).Value = "'140.52"
$ws.Cells.Item(29, 5).Value = '  -3.28%  '
$ws.Cells.Item(30, 4).Value = "'5.293"
$ws.Cells.Item(30, 5).Value = '  -7.99%  '
$ws.Cells.Item(31, 4).Value = "'7.860"
$ws.Cells.Item(31, 5).Value = '  -3.56%  '
$ws.Cells.Item(32, 4).Value = "'2.481"
$ws.Cells.Item(32, 5).Value = '  -0.67%  '
$ws.Cells.Item(33, 4).Value = "'1.801.56"
$ws.Cells.Item(33, 5).Value = '  -1.98%  '
$ws.Cells.Item(34, 4).Value = "'0.08140"
$ws.Cells.Item(34, 5).Value = '  +0.71%  '
$ws.Cells.Item(35, 4).Value = "'0.9911"
$ws.Cells.Item(35, 5).Value = '  -1.20%  '
$ws.Cells.Item(36, 4).Value = "'0.02905"
$ws.Cells.Item(36, 5).Value = '  -3.77%  '
$ws.Cells.Item(37, 4).Value = "'6.627"
$ws.Cells.Item(37, 5).Value = '  -3.37%  '
$ws.Cells.Item(38, 4).Value = "'0.2676"
$ws.Cells.Item(38, 5).Value = '  -3.03%  '
$ws.Cells.Item(39, 4).Value = "'0.09158"
$ws.Cells.Item(39, 5).Value = '  -3.29%  '
$ws.Cells.Item(40, 4).Value = "'10.34"
$ws.Cells.Item(40, 5).Value = '  +1.70%  '
$ws.Cells.Item(41, 4).Value = "'13.74"
$ws.Cells.Item(41, 5).Value = '  +3.26%  '
$ws.Cells.Item(42, 4).Value = "'1.432"
$ws.Cells.Item(42, 5).Value = '  -4.04%  '
$ws.Cells.Item(43, 4).Value = "'0.7543"
$ws.Cells.Item(43, 5).Value = '  -2.98%  '
$ws.Cells.Item(44, 4).Value = "'15.97"
$ws.Cells.Item(44, 5).Value = '  -1.34%  '
$ws.Cells.Item(45, 4).Value = "'0.6934"
$ws.Cells.Item(45, 5).Value = '  -0.98%  '
$ws.Cells.Item(46, 4).Value = "'2.477"
$ws.Cells.Item(46, 5).Value = '  -2.95%  '
$ws.Cells.Item(47, 4).Value = "'4.074"
$ws.Cells.Item(47, 5).Value = '  -1.58%  '
$ws.Cells.Item(48, 4).Value = "'0.9989"
$ws.Cells.Item(48, 5).Value = '  -0.45%  '
$ws.Cells.Item(49, 4).Value = "'0.08255"
$ws.Cells.Item(49, 5).Value = '  -3.42%  '
$ws.Cells.Item(50, 4).Value = "'135.85"
$ws.Cells.Item(50, 5).Value = '  -0.37%  '
$ws.Cells.Item(51, 4).Value = "'1.207"
$ws.Cells.Item(51, 5).Value = '  -6.97%  '
